# Update cryptocurrency price snapshot (column D) on Sheet1.
# Source values are stored as inline-string text (not numbers), e.g. "242.77",
# so we must write them back as text too. A bare numeric-looking string
# assigned via .Value would be auto-coerced to a number (and could pick up
# floating point noise, e.g. 23.14 -> 23.140000000000001), so we prefix with
# an apostrophe to force text entry, then reset the cell Style back to
# "Normal" so no stray number-format/quote-prefix styling is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue "D2"  "242.78"
Set-TextValue "D3"  "23.14"
Set-TextValue "D4"  "5.390"
Set-TextValue "D5"  "0.05979"
Set-TextValue "D6"  "3.403"
Set-TextValue "D7"  "6.466"
Set-TextValue "D8"  "0.8127"
Set-TextValue "D9"  "0.8998"
Set-TextValue "D10" "0.1414"
Set-TextValue "D11" "0.07413"
Set-TextValue "D12" "0.03356"
Set-TextValue "D13" "0.03071"
Set-TextValue "D14" "0.09341"
Set-TextValue "D15" "3.846"
Set-TextValue "D16" "0.001570"
Set-TextValue "D17" "0.04650"
Set-TextValue "D18" "0.0005941"
Set-TextValue "D19" "0.006097"
Set-TextValue "D20" "0.005019"
Set-TextValue "D21" "0.0009823"
Set-TextValue "D22" "0.00007800"
Set-TextValue "D23" "0.0002901"
Set-TextValue "D24" "3.615"
Set-TextValue "D25" "2.161"
Set-TextValue "D27" "0.1305"
Set-TextValue "D40" "0.03888"
Set-TextValue "D41" "0.006215"
Set-TextValue "D42" "0.1072"
Set-TextValue "D43" "0.002620"
Set-TextValue "D44" "0.007199"
Set-TextValue "D47" "0.0005801"
Set-TextValue "D49" "0.002298"
Set-TextValue "D50" "0.00002100"
Set-TextValue "D51" "0.0002000"
